$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Book" value for the first data row (D2): was YBB, now GG
$ws.Range("D2").Value = "GG"

# Clear the stray duplicate test rows (8-11), keep formatting, but drop values
$ws.Range("B8:G11").ClearContents()

# Move active selection back to D2
$ws.Range("D2").Select()
